# Mise à jour de l'application
# Adds 6 new "Entrainement" rows (J-1, 2025-09-19 / serial 45919) for the
# players who took part, appended after the existing last data row (521).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 521

$players = @(
    @{ Row = 522; Name = "Kamal Bafounta";   Poste = "center midfield"; Temps = "01:04:17";
       H = 4.84; I = 0.28; J = 4.56; K = 0.25; L = 0.04; M = 0; N = 0; O = 0;
       P = 4.5;  Q = 24.1;  R = 4.69; S = 24; T = 3;  U = 15; V = 4 },
    @{ Row = 523; Name = "Omar Benyounes";   Poste = "center midfield"; Temps = "01:04:24";
       H = 4.89; I = 0.22; J = 4.66; K = 0.21; L = 0.02; M = 0; N = 0; O = 0;
       P = 4.49; Q = 23.19; R = 3.9;  S = 21; T = 0;  U = 24; V = 1 },
    @{ Row = 524; Name = "Malik Boussaid";   Poste = "right back";      Temps = "00:55:51";
       H = 4.51; I = 0.26; J = 4.24; K = 0.26; L = 0.01; M = 0; N = 0; O = 0;
       P = 4.71; Q = 22.59; R = 4.61; S = 19; T = 4;  U = 19; V = 7 },
    @{ Row = 525; Name = "Mattheo Haon";     Poste = "right back";      Temps = "01:05:20";
       H = 5.11; I = 0.22; J = 4.88; K = 0.2;  L = 0.03; M = 0; N = 0; O = 0;
       P = 4.62; Q = 23.45; R = 4.61; S = 16; T = 4;  U = 27; V = 0 },
    @{ Row = 526; Name = "Ilan Ihaddadene";  Poste = "center midfield"; Temps = "01:05:14";
       H = 5.32; I = 0.14; J = 5.17; K = 0.12; L = 0.02; M = 0; N = 0; O = 1;
       P = 4.8;  Q = 25.71; R = 4.95; S = 18; T = 3;  U = 11; V = 1 },
    @{ Row = 527; Name = "Hedi Nasri";       Poste = "right back";      Temps = "01:04:29";
       H = 5.22; I = 0.29; J = 4.92; K = 0.22; L = 0.08; M = 0; N = 0; O = 1;
       P = 4.51; Q = 26.95; R = 5.15; S = 47; T = 11; U = 36; V = 11 }
)

foreach ($p in $players) {
    $r = $p.Row

    # Duplicate the formatting of the previous row (keeps B's date style,
    # D's center style, and the rest of the columns' default/no style).
    $ws.Range("A" + $lastRow + ":V" + $lastRow).Copy()
    $ws.Range("A" + $r + ":V" + $r).PasteSpecial(-4122)

    # Column E in these new rows carries style 6 (left/vcenter) - copy that
    # format specifically from an existing cell that already uses it.
    $ws.Range("E245").Copy()
    $ws.Range("E" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = "Entrainement"
    $ws.Cells.Item($r, 2).Value = 45919
    $ws.Cells.Item($r, 3).Value = "Global"
    $ws.Cells.Item($r, 4).Value = "J-1"
    $ws.Cells.Item($r, 5).Value = $p.Name
    $ws.Cells.Item($r, 6).Value = $p.Poste
    $ws.Cells.Item($r, 7).Value = $p.Temps
    $ws.Cells.Item($r, 8).Value = $p.H
    $ws.Cells.Item($r, 9).Value = $p.I
    $ws.Cells.Item($r, 10).Value = $p.J
    $ws.Cells.Item($r, 11).Value = $p.K
    $ws.Cells.Item($r, 12).Value = $p.L
    $ws.Cells.Item($r, 13).Value = $p.M
    $ws.Cells.Item($r, 14).Value = $p.N
    $ws.Cells.Item($r, 15).Value = $p.O
    $ws.Cells.Item($r, 16).Value = $p.P
    $ws.Cells.Item($r, 17).Value = $p.Q
    $ws.Cells.Item($r, 18).Value = $p.R
    $ws.Cells.Item($r, 19).Value = $p.S
    $ws.Cells.Item($r, 20).Value = $p.T
    $ws.Cells.Item($r, 21).Value = $p.U
    $ws.Cells.Item($r, 22).Value = $p.V
}

# Restore selection/cursor to match where Excel would land after the paste
# (one row below the new last data row, column E).
$ws.Range("E532").Select()
